$wb = $excel.ActiveWorkbook

# --- "Data" sheet: append the newest WRESBAL observation ---
$data = $wb.Worksheets.Item("Data")
$data.Cells.Item(112, 1).Value = Get-Date -Year 2023 -Month 11 -Day 15
$data.Cells.Item(112, 2).Value = 3391.713

# --- "SeriesInfo" sheet: refresh the metadata pulled from FRED ---
$info = $wb.Worksheets.Item("SeriesInfo")
$info.Range("B3").Value = "2023-11-21"
$info.Range("B4").Value = "2023-11-21"
$info.Range("B7").Value = "2023-11-15"
$info.Range("B14").Value = "2023-11-16 15:36:02-06"
$info.Range("B15").Value = 74
